# feat: add 2022-Q4 data
#
# - Inserts a new "2022-Q4" sheet between "总计" and "2022-Q3", populated
#   with the Q4 fund-holdings table.
# - Updates the "总计" (totals) sheet: the former single "2022-Q3" summary
#   row becomes row 3, and a new row 2 is added for "2022-Q4" totals
#   (14 funds held, 2.58 亿元 market value).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsQ3    = $wb.Worksheets.Item(2)   # "2022-Q3" (currently 2nd sheet)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right before "2022-Q3" so the final
#    tab order is: 总计, 2022-Q4, 2022-Q3
# ---------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Add($wsQ3)
$wsQ4.Name = "2022-Q4"

# Header row (B1:H1) - bold, centered, bordered, matching the workbook's
# existing header style.
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$col = 2
foreach ($h in $headers) {
    $c = $wsQ4.Cells.Item(1, $col)
    $c.Value = $h
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.LineStyle = 1
    $col++
}

# Fund holding rows (2022-Q4 data). Columns D/E/F/G are stored as plain
# text (leading apostrophe forces text so values like "24.00" / "012526"
# keep their exact formatting instead of being coerced to numbers).
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").Value = "'012526"
$wsQ4.Range("C2").Value = "广发盛锦混合A"
$wsQ4.Range("D2").Value = "'24.00"
$wsQ4.Range("E2").Value = "'93.13"
$wsQ4.Range("F2").Value = "'4.90"
$wsQ4.Range("G2").Value = "'1.1760"
$wsQ4.Range("H2").Value = 5
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").Value = "'013141"
$wsQ4.Range("C3").Value = "中信保诚弘远混合A"
$wsQ4.Range("D3").Value = "'18.60"
$wsQ4.Range("E3").Value = "'70.64"
$wsQ4.Range("F3").Value = "'2.98"
$wsQ4.Range("G3").Value = "'0.5543"
$wsQ4.Range("H3").Value = 3
$wsQ4.Range("A4").Value = 2
$wsQ4.Range("B4").Value = "'165516"
$wsQ4.Range("C4").Value = "信诚周期轮动混合（LOF）A"
$wsQ4.Range("D4").Value = "'17.18"
$wsQ4.Range("E4").Value = "'70.42"
$wsQ4.Range("F4").Value = "'2.96"
$wsQ4.Range("G4").Value = "'0.5085"
$wsQ4.Range("H4").Value = 4
$wsQ4.Range("A5").Value = 3
$wsQ4.Range("B5").Value = "'005682"
$wsQ4.Range("C5").Value = "财通资管消费精选灵活配置混合A"
$wsQ4.Range("D5").Value = "'3.44"
$wsQ4.Range("E5").Value = "'95.15"
$wsQ4.Range("F5").Value = "'4.06"
$wsQ4.Range("G5").Value = "'0.1397"
$wsQ4.Range("H5").Value = 10
$wsQ4.Range("A6").Value = 4
$wsQ4.Range("B6").Value = "'012527"
$wsQ4.Range("C6").Value = "广发盛锦混合C"
$wsQ4.Range("D6").Value = "'1.14"
$wsQ4.Range("E6").Value = "'93.13"
$wsQ4.Range("F6").Value = "'4.90"
$wsQ4.Range("G6").Value = "'0.0559"
$wsQ4.Range("H6").Value = 5
$wsQ4.Range("A7").Value = 5
$wsQ4.Range("B7").Value = "'015769"
$wsQ4.Range("C7").Value = "天弘低碳经济混合A"
$wsQ4.Range("D7").Value = "'1.12"
$wsQ4.Range("E7").Value = "'86.07"
$wsQ4.Range("F7").Value = "'3.40"
$wsQ4.Range("G7").Value = "'0.0381"
$wsQ4.Range("H7").Value = 10
$wsQ4.Range("A8").Value = 6
$wsQ4.Range("B8").Value = "'008277"
$wsQ4.Range("C8").Value = "财通资管行业精选混合"
$wsQ4.Range("D8").Value = "'0.95"
$wsQ4.Range("E8").Value = "'90.07"
$wsQ4.Range("F8").Value = "'3.88"
$wsQ4.Range("G8").Value = "'0.0369"
$wsQ4.Range("H8").Value = 10
$wsQ4.Range("A9").Value = 7
$wsQ4.Range("B9").Value = "'015770"
$wsQ4.Range("C9").Value = "天弘低碳经济混合C"
$wsQ4.Range("D9").Value = "'0.99"
$wsQ4.Range("E9").Value = "'86.07"
$wsQ4.Range("F9").Value = "'3.40"
$wsQ4.Range("G9").Value = "'0.0337"
$wsQ4.Range("H9").Value = 10
$wsQ4.Range("A10").Value = 8
$wsQ4.Range("B10").Value = "'001261"
$wsQ4.Range("C10").Value = "中融新机遇灵活配置混合"
$wsQ4.Range("D10").Value = "'0.46"
$wsQ4.Range("E10").Value = "'78.82"
$wsQ4.Range("F10").Value = "'3.87"
$wsQ4.Range("G10").Value = "'0.0178"
$wsQ4.Range("H10").Value = 6
$wsQ4.Range("A11").Value = 9
$wsQ4.Range("B11").Value = "'011020"
$wsQ4.Range("C11").Value = "财通资管消费精选灵活配置混合C"
$wsQ4.Range("D11").Value = "'0.16"
$wsQ4.Range("E11").Value = "'95.15"
$wsQ4.Range("F11").Value = "'4.06"
$wsQ4.Range("G11").Value = "'0.0065"
$wsQ4.Range("H11").Value = 10
$wsQ4.Range("A12").Value = 10
$wsQ4.Range("B12").Value = "'012287"
$wsQ4.Range("C12").Value = "东海启航6个月持有期混合A"
$wsQ4.Range("D12").Value = "'0.74"
$wsQ4.Range("E12").Value = "'28.90"
$wsQ4.Range("F12").Value = "'0.65"
$wsQ4.Range("G12").Value = "'0.0048"
$wsQ4.Range("H12").Value = 10
$wsQ4.Range("A13").Value = 11
$wsQ4.Range("B13").Value = "'013377"
$wsQ4.Range("C13").Value = "东海启航6个月持有期混合C"
$wsQ4.Range("D13").Value = "'0.23"
$wsQ4.Range("E13").Value = "'28.90"
$wsQ4.Range("F13").Value = "'0.65"
$wsQ4.Range("G13").Value = "'0.0015"
$wsQ4.Range("H13").Value = 10
$wsQ4.Range("A14").Value = 12
$wsQ4.Range("B14").Value = "'014335"
$wsQ4.Range("C14").Value = "信诚周期轮动混合（LOF）C"
$wsQ4.Range("D14").Value = "'0.04"
$wsQ4.Range("E14").Value = "'70.42"
$wsQ4.Range("F14").Value = "'2.96"
$wsQ4.Range("G14").Value = "'0.0012"
$wsQ4.Range("H14").Value = 4
$wsQ4.Range("A15").Value = 13
$wsQ4.Range("B15").Value = "'015936"
$wsQ4.Range("C15").Value = "中信保诚弘远混合C"
$wsQ4.Range("D15").Value = "'0.02"
$wsQ4.Range("E15").Value = "'70.64"
$wsQ4.Range("F15").Value = "'2.98"
$wsQ4.Range("G15").Value = "'0.0006"
$wsQ4.Range("H15").Value = 3

# Style column A (row index) the same as the header cells.
for ($r = 2; $r -le 15; $r++) {
    $c = $wsQ4.Cells.Item($r, 1)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing "2022-Q3" row down to
#    row 3, keeping its style, then write the new "2022-Q4" row 2.
# ---------------------------------------------------------------------
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 5
$wsTotal.Range("D3").Value = 0.86

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 14
$wsTotal.Range("D2").Value = 2.58
